$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44377   # D2
$ws.Cells.Item(2, 10).Value = 80  # J2
$ws.Cells.Item(2, 11).Value = 18000  # K2
$ws.Cells.Item(2, 12).Value = 19000  # L2
$ws.Cells.Item(2, 13).Value = 18500  # M2
$ws.Cells.Item(2, 16).Value = 1233  # P2

$ws.Cells.Item(3, 4).Value = 44326   # D3
$ws.Cells.Item(3, 10).Value = 45  # J3
$ws.Cells.Item(3, 11).Value = 15000  # K3
$ws.Cells.Item(3, 12).Value = 15000  # L3
$ws.Cells.Item(3, 13).Value = 15000  # M3
$ws.Cells.Item(3, 16).Value = 1000  # P3

$ws.Cells.Item(4, 4).Value = 44314   # D4
$ws.Cells.Item(4, 10).Value = 45  # J4
$ws.Cells.Item(4, 11).Value = 15000  # K4
$ws.Cells.Item(4, 12).Value = 15000  # L4
$ws.Cells.Item(4, 13).Value = 15000  # M4
$ws.Cells.Item(4, 16).Value = 1000  # P4

$ws.Cells.Item(5, 4).Value = 44323   # D5
$ws.Cells.Item(5, 10).Value = 40  # J5
$ws.Cells.Item(5, 11).Value = 15000  # K5
$ws.Cells.Item(5, 12).Value = 15000  # L5
$ws.Cells.Item(5, 13).Value = 15000  # M5
$ws.Cells.Item(5, 16).Value = 1000  # P5

$ws.Cells.Item(6, 4).Value = 44370   # D6
$ws.Cells.Item(6, 10).Value = 50  # J6
$ws.Cells.Item(6, 11).Value = 18000  # K6
$ws.Cells.Item(6, 12).Value = 18000  # L6
$ws.Cells.Item(6, 13).Value = 18000  # M6
$ws.Cells.Item(6, 16).Value = 1200  # P6

$ws.Cells.Item(7, 4).Value = 44315   # D7
$ws.Cells.Item(7, 10).Value = 65  # J7
$ws.Cells.Item(7, 11).Value = 14000  # K7
$ws.Cells.Item(7, 12).Value = 15000  # L7
$ws.Cells.Item(7, 13).Value = 14538  # M7
$ws.Cells.Item(7, 16).Value = 969  # P7

$ws.Cells.Item(8, 4).Value = 44406   # D8
$ws.Cells.Item(8, 10).Value = 50  # J8
$ws.Cells.Item(8, 11).Value = 22000  # K8
$ws.Cells.Item(8, 12).Value = 22000  # L8
$ws.Cells.Item(8, 13).Value = 22000  # M8
$ws.Cells.Item(8, 16).Value = 1467  # P8

$ws.Cells.Item(9, 4).Value = 44344   # D9
$ws.Cells.Item(9, 10).Value = 40  # J9
$ws.Cells.Item(9, 11).Value = 20000  # K9
$ws.Cells.Item(9, 12).Value = 20000  # L9
$ws.Cells.Item(9, 13).Value = 20000  # M9
$ws.Cells.Item(9, 16).Value = 1333  # P9

$ws.Cells.Item(10, 4).Value = 44399   # D10
$ws.Cells.Item(10, 10).Value = 38  # J10
$ws.Cells.Item(10, 11).Value = 22000  # K10
$ws.Cells.Item(10, 12).Value = 22000  # L10
$ws.Cells.Item(10, 13).Value = 22000  # M10
$ws.Cells.Item(10, 16).Value = 1467  # P10

$ws.Cells.Item(11, 4).Value = 44333   # D11
$ws.Cells.Item(11, 10).Value = 35  # J11
$ws.Cells.Item(11, 11).Value = 15000  # K11
$ws.Cells.Item(11, 12).Value = 15000  # L11
$ws.Cells.Item(11, 13).Value = 15000  # M11
$ws.Cells.Item(11, 16).Value = 1000  # P11

$ws.Cells.Item(12, 4).Value = 44340   # D12
$ws.Cells.Item(12, 10).Value = 47  # J12
$ws.Cells.Item(12, 11).Value = 14000  # K12
$ws.Cells.Item(12, 12).Value = 14000  # L12
$ws.Cells.Item(12, 13).Value = 14000  # M12
$ws.Cells.Item(12, 16).Value = 933  # P12

$ws.Cells.Item(13, 4).Value = 44321   # D13
$ws.Cells.Item(13, 10).Value = 38  # J13
$ws.Cells.Item(13, 11).Value = 15000  # K13
$ws.Cells.Item(13, 12).Value = 15000  # L13
$ws.Cells.Item(13, 13).Value = 15000  # M13
$ws.Cells.Item(13, 16).Value = 1000  # P13

$ws.Cells.Item(14, 4).Value = 44455   # D14
$ws.Cells.Item(14, 10).Value = 35  # J14
$ws.Cells.Item(14, 11).Value = 22000  # K14
$ws.Cells.Item(14, 12).Value = 22000  # L14
$ws.Cells.Item(14, 13).Value = 22000  # M14
$ws.Cells.Item(14, 16).Value = 1467  # P14

$ws.Cells.Item(15, 4).Value = 44397   # D15
$ws.Cells.Item(15, 10).Value = 73  # J15
$ws.Cells.Item(15, 11).Value = 21000  # K15
$ws.Cells.Item(15, 12).Value = 22000  # L15
$ws.Cells.Item(15, 13).Value = 21521  # M15
$ws.Cells.Item(15, 16).Value = 1435  # P15

$ws.Cells.Item(16, 4).Value = 44336   # D16
$ws.Cells.Item(16, 10).Value = 65  # J16
$ws.Cells.Item(16, 11).Value = 14000  # K16
$ws.Cells.Item(16, 12).Value = 15000  # L16
$ws.Cells.Item(16, 13).Value = 14462  # M16
$ws.Cells.Item(16, 16).Value = 964  # P16

$ws.Cells.Item(17, 4).Value = 44309   # D17
$ws.Cells.Item(17, 10).Value = 50  # J17
$ws.Cells.Item(17, 11).Value = 15000  # K17
$ws.Cells.Item(17, 12).Value = 15000  # L17
$ws.Cells.Item(17, 13).Value = 15000  # M17
$ws.Cells.Item(17, 16).Value = 1000  # P17

$ws.Cells.Item(18, 4).Value = 44319   # D18
$ws.Cells.Item(18, 10).Value = 50  # J18
$ws.Cells.Item(18, 11).Value = 15000  # K18
$ws.Cells.Item(18, 12).Value = 15000  # L18
$ws.Cells.Item(18, 13).Value = 15000  # M18
$ws.Cells.Item(18, 16).Value = 1000  # P18

$ws.Cells.Item(19, 4).Value = 44438   # D19
$ws.Cells.Item(19, 10).Value = 75  # J19
$ws.Cells.Item(19, 11).Value = 19000  # K19
$ws.Cells.Item(19, 12).Value = 20000  # L19
$ws.Cells.Item(19, 13).Value = 19467  # M19
$ws.Cells.Item(19, 16).Value = 1298  # P19

$ws.Cells.Item(20, 4).Value = 44308   # D20
$ws.Cells.Item(20, 10).Value = 40  # J20
$ws.Cells.Item(20, 11).Value = 16000  # K20
$ws.Cells.Item(20, 12).Value = 16000  # L20
$ws.Cells.Item(20, 13).Value = 16000  # M20
$ws.Cells.Item(20, 16).Value = 1067  # P20

$ws.Cells.Item(21, 4).Value = 44320   # D21
$ws.Cells.Item(21, 10).Value = 40  # J21
$ws.Cells.Item(21, 11).Value = 15000  # K21
$ws.Cells.Item(21, 12).Value = 15000  # L21
$ws.Cells.Item(21, 13).Value = 15000  # M21
$ws.Cells.Item(21, 16).Value = 1000  # P21

$ws.Cells.Item(22, 4).Value = 44343   # D22
$ws.Cells.Item(22, 10).Value = 40  # J22
$ws.Cells.Item(22, 11).Value = 15000  # K22
$ws.Cells.Item(22, 12).Value = 15000  # L22
$ws.Cells.Item(22, 13).Value = 15000  # M22
$ws.Cells.Item(22, 16).Value = 1000  # P22

$ws.Cells.Item(23, 4).Value = 44316   # D23
$ws.Cells.Item(23, 10).Value = 45  # J23
$ws.Cells.Item(23, 11).Value = 14000  # K23
$ws.Cells.Item(23, 12).Value = 15000  # L23
$ws.Cells.Item(23, 13).Value = 14444  # M23
$ws.Cells.Item(23, 16).Value = 963  # P23

$ws.Cells.Item(24, 4).Value = 44329   # D24
$ws.Cells.Item(24, 10).Value = 35  # J24
$ws.Cells.Item(24, 11).Value = 15000  # K24
$ws.Cells.Item(24, 12).Value = 15000  # L24
$ws.Cells.Item(24, 13).Value = 15000  # M24
$ws.Cells.Item(24, 16).Value = 1000  # P24

$ws.Cells.Item(25, 4).Value = 44411   # D25
$ws.Cells.Item(25, 10).Value = 50  # J25
$ws.Cells.Item(25, 11).Value = 22000  # K25
$ws.Cells.Item(25, 12).Value = 22000  # L25
$ws.Cells.Item(25, 13).Value = 22000  # M25
$ws.Cells.Item(25, 16).Value = 1467  # P25

$ws.Cells.Item(26, 4).Value = 44313   # D26
$ws.Cells.Item(26, 10).Value = 40  # J26
$ws.Cells.Item(26, 11).Value = 14000  # K26
$ws.Cells.Item(26, 12).Value = 14000  # L26
$ws.Cells.Item(26, 13).Value = 14000  # M26
$ws.Cells.Item(26, 16).Value = 933  # P26

$ws.Cells.Item(27, 4).Value = 44334   # D27
$ws.Cells.Item(27, 10).Value = 50  # J27
$ws.Cells.Item(27, 11).Value = 14000  # K27
$ws.Cells.Item(27, 12).Value = 14000  # L27
$ws.Cells.Item(27, 13).Value = 14000  # M27
$ws.Cells.Item(27, 16).Value = 933  # P27

$ws.Cells.Item(28, 4).Value = 44330   # D28
$ws.Cells.Item(28, 10).Value = 30  # J28
$ws.Cells.Item(28, 11).Value = 15000  # K28
$ws.Cells.Item(28, 12).Value = 15000  # L28
$ws.Cells.Item(28, 13).Value = 15000  # M28
$ws.Cells.Item(28, 16).Value = 1000  # P28

$ws.Cells.Item(29, 4).Value = 44448   # D29
$ws.Cells.Item(29, 10).Value = 85  # J29
$ws.Cells.Item(29, 11).Value = 21000  # K29
$ws.Cells.Item(29, 12).Value = 22000  # L29
$ws.Cells.Item(29, 13).Value = 21529  # M29
$ws.Cells.Item(29, 16).Value = 1435  # P29

$ws.Cells.Item(30, 4).Value = 44328   # D30
$ws.Cells.Item(30, 10).Value = 38  # J30
$ws.Cells.Item(30, 11).Value = 15000  # K30
$ws.Cells.Item(30, 12).Value = 15000  # L30
$ws.Cells.Item(30, 13).Value = 15000  # M30
$ws.Cells.Item(30, 16).Value = 1000  # P30

$ws.Cells.Item(31, 4).Value = 44341   # D31
$ws.Cells.Item(31, 10).Value = 40  # J31
$ws.Cells.Item(31, 11).Value = 15000  # K31
$ws.Cells.Item(31, 12).Value = 15000  # L31
$ws.Cells.Item(31, 13).Value = 15000  # M31
$ws.Cells.Item(31, 16).Value = 1000  # P31

$ws.Cells.Item(32, 4).Value = 44452   # D32
$ws.Cells.Item(32, 10).Value = 73  # J32
$ws.Cells.Item(32, 11).Value = 22000  # K32
$ws.Cells.Item(32, 12).Value = 23000  # L32
$ws.Cells.Item(32, 13).Value = 22479  # M32
$ws.Cells.Item(32, 16).Value = 1499  # P32

$ws.Cells.Item(33, 4).Value = 44322   # D33
$ws.Cells.Item(33, 10).Value = 70  # J33
$ws.Cells.Item(33, 11).Value = 14000  # K33
$ws.Cells.Item(33, 12).Value = 15000  # L33
$ws.Cells.Item(33, 13).Value = 14500  # M33
$ws.Cells.Item(33, 16).Value = 967  # P33

$ws.Cells.Item(34, 4).Value = 44327   # D34
$ws.Cells.Item(34, 10).Value = 35  # J34
$ws.Cells.Item(34, 11).Value = 15000  # K34
$ws.Cells.Item(34, 12).Value = 15000  # L34
$ws.Cells.Item(34, 13).Value = 15000  # M34
$ws.Cells.Item(34, 16).Value = 1000  # P34

$ws.Cells.Item(35, 4).Value = 44312   # D35
$ws.Cells.Item(35, 10).Value = 80  # J35
$ws.Cells.Item(35, 11).Value = 13000  # K35
$ws.Cells.Item(35, 12).Value = 14000  # L35
$ws.Cells.Item(35, 13).Value = 13562  # M35
$ws.Cells.Item(35, 16).Value = 904  # P35
